# "alguns testes na planilha"
$wb = $excel.ActiveWorkbook

# --- alunos: add two new student rows (8 and 9) ---
$ws1 = $wb.Worksheets.Item("alunos")

# Make sure the new cells don't inherit the column styles (rows 5-7 above
# them are also style-less), so reset to Normal before writing values.
$ws1.Range("A8:G9").Style = "Normal"

$ws1.Range("A8").Value = 7
$ws1.Range("B8").Value = "Wellington almeida"
$ws1.Range("C8").Value = "7 ani"
$ws1.Range("D8").Value = "Tarde"
$ws1.Range("E8").Value = 20
$ws1.Range("F8").Value = "88 9 81762299"
$ws1.Range("G8").Value = "Rua 1"

$ws1.Range("A9").Value = 8
$ws1.Range("B9").Value = "Rian"
$ws1.Range("C9").Value = "9 ano"
$ws1.Range("D9").Value = "Tarde"
$ws1.Range("E9").Value = 20
$ws1.Range("F9").Value = "88 99 99 99999"
$ws1.Range("G9").Value = "Rua 3"

# --- emprestimos: remove the test/extra loan entry in row 3, keep the row ---
$ws3 = $wb.Worksheets.Item("emprestimos")
$ws3.Range("A3:D3").ClearContents()
